$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.786.88'
$ws.Range('E2').Value = '  -0.62%  '
$ws.Range('D3').Value = '2.454.76'
$ws.Range('E3').Value = '  -1.32%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '557.80'
$ws.Range('E5').Value = '  -1.72%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '162.20'
$ws.Range('E6').Value = '  -2.10%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.504'
$ws.Range('E8').Value = '  -1.18%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.152'
$ws.Range('E9').Value = '  -4.02%  '
$ws.Range('D10').Value = '2.288.77'
$ws.Range('E10').Value = '  -7.99%  '
$ws.Range('E11').Value = '  -0.43%  '
$ws.Range('E12').Value = '  -3.52%  '
$ws.Range('E13').Value = '  -0.93%  '
$ws.Range('D14').Value = '2.905.97'
$ws.Range('E14').Value = '  -1.28%  '
$ws.Range('D15').Value = '68.626.98'
$ws.Range('E15').Value = '  -0.72%  '
$ws.Range('E16').Value = '  -2.87%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '23.55'
$ws.Range('E17').Value = '  -2.00%  '
$ws.Range('D18').Value = '2.472.82'
$ws.Range('E18').Value = '  -0.97%  '
$ws.Range('E19').Value = '  -3.59%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '341.17'
$ws.Range('E20').Value = '  -3.12%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.06'
$ws.Range('E21').Value = '  -3.83%  '
$ws.Range('E22').Value = '  -2.48%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.93'
$ws.Range('E23').Value = '  +0.94%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.00'
$ws.Range('E24').Value = '  +0.06%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '66.93'
$ws.Range('E25').Value = '  -3.09%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.69'
$ws.Range('E26').Value = '  -2.58%  '
$ws.Range('D27').Value = '2.581.19'
$ws.Range('E27').Value = '  -1.35%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.993'
$ws.Range('E28').Value = '  -1.23%  '
$ws.Range('E29').Value = '  -4.84%  '
$ws.Range('D30').Value = '0.0₃0816'
$ws.Range('E30').Value = '  -5.91%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.16'
$ws.Range('E31').Value = '  -4.14%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '435.39'
$ws.Range('E32').Value = '  -0.45%  '
$ws.Range('E33').Value = '  +0.00%  '
$ws.Range('E34').Value = '  -3.61%  '
$ws.Range('E35').Value = '  -5.64%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '157.12'
$ws.Range('E36').Value = '  +1.71%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '19.03'
$ws.Range('E37').Value = '  -0.04%  '
$ws.Range('E38').Value = '  +0.09%  '
$ws.Range('E39').Value = '  -3.57%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '17.81'
$ws.Range('E40').Value = '  -1.52%  '
$ws.Range('E41').Value = '  -2.49%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.301'
$ws.Range('E42').Value = '  -3.53%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '37.43'
$ws.Range('E43').Value = '  -1.06%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.48'
$ws.Range('E44').Value = '  -5.60%  '
$ws.Range('E45').Value = '  +3.49%  '
$ws.Range('E46').Value = '  -4.49%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '132.97'
$ws.Range('E47').Value = '  -3.73%  '
$ws.Range('E48').Value = '  -2.17%  '
$ws.Range('E49').Value = '  -0.63%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.482'
$ws.Range('E50').Value = '  -4.23%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.559'
$ws.Range('E51').Value = '  -2.40%  '
